# Swap the order of the "System" / email entries in column G
# (Recorded By) wherever the cell currently reads
# "System, dnasr281@gmail.com", changing it to
# "dnasr281@gmail.com, System".
#
# Using Find/FindNext (rather than iterating every cell in the sheet)
# avoids touching the many blank-but-styled cells that also live in
# column G, so only the cells that actually contain the old text get
# written.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$first = $ws.Cells.Find($oldValue)

if ($first -ne $null) {
    $visited = @()
    $cur = $first
    $addr = $cur.Row.ToString() + ":" + $cur.Column.ToString()
    $guard = 0

    while (($visited -notcontains $addr) -and ($guard -lt 1000)) {
        $visited += $addr
        $guard++

        $cur.Value2 = $newValue

        $cur = $ws.Cells.FindNext($cur)
        if ($cur -eq $null) { break }
        $addr = $cur.Row.ToString() + ":" + $cur.Column.ToString()
    }
}
